$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current last data row (row 19: "99" / "99999" for
# "gebied onbekend"), shifting it down to become row 20.
$ws.Rows.Item(19).Insert()

# Fill the newly inserted row 19 with the "uitrustingsgraad/niveau" category
# (logo 93 / provincie-code 99993), matching the number formatting of the
# surrounding data rows.
$ws.Cells.Item(19, 1).Value = 93
$ws.Cells.Item(19, 2).Value = 99993
$ws.Cells.Item(19, 1).NumberFormat = $ws.Cells.Item(18, 1).NumberFormat
$ws.Cells.Item(19, 2).NumberFormat = $ws.Cells.Item(18, 2).NumberFormat
